$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing 7 detail rows (old rows 17-23) that listed the
# periods/amounts for the other employees. This shifts the signature
# block (old rows 28/29) up to rows 21/22.
$ws.Rows("17:23").Delete()

# The surviving data row (row 16) now carries the figures that used to
# belong to the last remaining worker's record (previously row 23).
$ws.Range("F16").Value = 28000
$ws.Range("G16").Value = 1000000

# Update the summary/header figures to reflect the now-single employee.
$ws.Range("E11").Value = 28000
$ws.Range("C13").Value = 1
$ws.Range("F13").Value = 1

# Column D ("Nombre Trabajador") no longer needs to fit the longest of
# the old names, so its best-fit width shrinks.
$ws.Columns("D").ColumnWidth = 25.46
